$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (05-nov) before the
#     existing "01-oct." column (DJ), shifting DJ:EN right to DK:EO. ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")
$wsSpot.Columns("DJ").Insert()
$wsSpot.Range("DJ1").Value = "05-nov"
$wsSpot.Range("DJ2:DJ25").Value = "-"

# --- Sheet "Gaz": append the new day's row. ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A142").Value = "'2025-11-03"
$wsGaz.Range("B142").Value = 30.35

# --- Sheet "CO2": append the new day's row. ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A142").Value = "'2025-11-03"
$wsCo2.Range("B142").Value = 81.2

Write-Output "edit applied"
